$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.947.64"
$ws.Range("E2").Value = "  -3.58%  "
$ws.Range("D3").Value = "3.527.18"
$ws.Range("E3").Value = "  -3.87%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "610.33"
$ws.Range("E5").Value = "  -5.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.10"
$ws.Range("E6").Value = "  -4.25%  "
$ws.Range("D7").Value = "3.525.81"
$ws.Range("E7").Value = "  -3.85%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  -3.12%  "
$ws.Range("E10").Value = "  -3.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.84"
$ws.Range("E11").Value = "  -3.63%  "
$ws.Range("E12").Value = "  -3.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000219"
$ws.Range("E13").Value = "  -4.81%  "
$ws.Range("D14").Value = "4.127.43"
$ws.Range("E14").Value = "  -3.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.63"
$ws.Range("E15").Value = "  -2.68%  "
$ws.Range("D16").Value = "3.518.90"
$ws.Range("E16").Value = "  -4.22%  "
$ws.Range("D17").Value = "66.935.63"
$ws.Range("E17").Value = "  -3.55%  "
$ws.Range("E18").Value = "  +0.60%  "
$ws.Range("E19").Value = "  -2.55%  "
$ws.Range("E20").Value = "  -3.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "443.40"
$ws.Range("E21").Value = "  -5.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.19"
$ws.Range("E22").Value = "  -8.37%  "
$ws.Range("E23").Value = "  -2.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.71"
$ws.Range("E24").Value = "  -1.98%  "
$ws.Range("D25").Value = "3.670.11"
$ws.Range("E25").Value = "  -3.79%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000122"
$ws.Range("E27").Value = "  -1.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.21"
$ws.Range("E28").Value = "  -6.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.17"
$ws.Range("E29").Value = "  -9.08%  "
$ws.Range("E30").Value = "  -3.50%  "
$ws.Range("E31").Value = "  -1.28%  "
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "25.71"
$ws.Range("E33").Value = "  -3.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.159"
$ws.Range("E34").Value = "  -2.52%  "
$ws.Range("E35").Value = "  -4.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.86"
$ws.Range("E36").Value = "  -6.19%  "
$ws.Range("D37").Value = "3.523.27"
$ws.Range("E37").Value = "  -3.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.00"
$ws.Range("E38").Value = "  -4.78%  "
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "174.20"
$ws.Range("E41").Value = "  -2.42%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.13"
$ws.Range("E42").Value = "  -2.28%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.56"
$ws.Range("E43").Value = "  -5.06%  "
$ws.Range("E44").Value = "  -3.75%  "
$ws.Range("E45").Value = "  -3.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.32"
$ws.Range("E46").Value = "  -3.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.18"
$ws.Range("E47").Value = "  -5.53%  "
$ws.Range("E48").Value = "  -4.68%  "
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.56"
$ws.Range("E50").Value = "  -2.80%  "
$ws.Range("E51").Value = "  -3.43%  "
